$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# row 13
$ws.Range("H13").Value2 = 8032.364
$ws.Range("I13").Value2 = 17950
$ws.Range("J13").Value2 = 5828.4443
$ws.Range("K13").Value2 = 17950
$ws.Range("L13").Value2 = 5828.4443
$ws.Range("M13").Value2 = -17781
$ws.Range("N13").Value2 = -6166.4443
# row 43
$ws.Range("H43").Value2 = 2776.5557
$ws.Range("I43").Value2 = 3638
$ws.Range("J43").Value2 = 1699.75
$ws.Range("K43").Value2 = 3638
$ws.Range("L43").Value2 = 1699.75
$ws.Range("M43").Value2 = -3569
$ws.Range("N43").Value2 = -1837.75
# row 54
$ws.Range("H54").Value2 = 2858.6667
$ws.Range("I54").Value2 = 2858.6667
$ws.Range("J54").Value2 = 0
$ws.Range("K54").Value2 = 2858.6667
$ws.Range("L54").Value2 = 0
$ws.Range("M54").Value2 = -2372.6667
$ws.Range("N54").ClearContents()
# row 97
$ws.Range("H97").Value2 = 1637.1666
$ws.Range("I97").Value2 = 950
$ws.Range("J97").Value2 = 1980.75
$ws.Range("K97").Value2 = 2850
$ws.Range("L97").Value2 = 5942.25
$ws.Range("M97").Value2 = -2354
$ws.Range("N97").Value2 = -6934.25
# row 112
$ws.Range("H112").Value2 = 1670.8695
$ws.Range("J112").Value2 = 1766
$ws.Range("L112").Value2 = 5298
$ws.Range("N112").Value2 = -7514
# row 137
$ws.Range("H137").Value2 = 682.4286
$ws.Range("J137").Value2 = 767.6667
$ws.Range("L137").Value2 = 2303.0001
$ws.Range("N137").Value2 = -7403.0001

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value2 = 4211.78
$ws.Range("I32").Value2 = 3523.1702
$ws.Range("J32").Value2 = 15000
$ws.Range("K32").Value2 = 3523.1702
$ws.Range("L32").Value2 = 15000
$ws.Range("M32").Value2 = -3236.1702
$ws.Range("N32").Value2 = -15574
# row 109
$ws.Range("H109").Value2 = 43325
$ws.Range("J109").Value2 = 43325
$ws.Range("L109").Value2 = 43325
$ws.Range("N109").Value2 = -46099

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# row 7
$ws.Range("H7").Value2 = 101.5
$ws.Range("I7").Value2 = 101.5
$ws.Range("K7").Value2 = 101.5
$ws.Range("M7").Value2 = 11.5
# row 20
$ws.Range("H20").Value2 = 48099.31
$ws.Range("I20").Value2 = 121398.2
$ws.Range("J20").Value2 = 2287.5
$ws.Range("K20").Value2 = 121398.2
$ws.Range("L20").Value2 = 2287.5
$ws.Range("M20").Value2 = -121151.2
$ws.Range("N20").Value2 = -2781.5
# row 107
$ws.Range("H107").Value2 = 16550.666
$ws.Range("I107").Value2 = 1386.8
$ws.Range("J107").Value2 = 35505.5
$ws.Range("K107").Value2 = 1386.8
$ws.Range("L107").Value2 = 35505.5
$ws.Range("M107").Value2 = 533.2
$ws.Range("N107").Value2 = -39345.5

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# row 2
$ws.Range("H2").Value2 = 225
$ws.Range("J2").Value2 = 0
$ws.Range("L2").Value2 = 0
$ws.Range("N2").ClearContents()
# row 4
$ws.Range("H4").Value2 = 1998.1818
$ws.Range("J4").Value2 = 1998.1818
$ws.Range("L4").Value2 = 1998.1818
$ws.Range("N4").Value2 = -2222.1818
# row 31
$ws.Range("H31").Value2 = 1798.8718
$ws.Range("I31").Value2 = 1746.2106
$ws.Range("J31").Value2 = 3800
$ws.Range("K31").Value2 = 1746.2106
$ws.Range("L31").Value2 = 3800
$ws.Range("M31").Value2 = -1451.2106
$ws.Range("N31").Value2 = -4390
# row 34
$ws.Range("H34").Value2 = 1798.8718
$ws.Range("I34").Value2 = 1746.2106
$ws.Range("J34").Value2 = 3800
$ws.Range("K34").Value2 = 1746.2106
$ws.Range("L34").Value2 = 3800
$ws.Range("M34").Value2 = -1544.2106
$ws.Range("N34").Value2 = -4204
# row 42
$ws.Range("H42").Value2 = 7500
$ws.Range("J42").Value2 = 10000
$ws.Range("L42").Value2 = 10000
$ws.Range("N42").Value2 = -11186
# row 99
$ws.Range("H99").Value2 = 2128.5715
$ws.Range("I99").Value2 = 2040
$ws.Range("J99").Value2 = 2350
$ws.Range("K99").Value2 = 2040
$ws.Range("L99").Value2 = 2350
$ws.Range("M99").Value2 = -542
$ws.Range("N99").Value2 = -5346
# row 103
$ws.Range("H103").Value2 = 4000
$ws.Range("I103").Value2 = 4000
$ws.Range("K103").Value2 = 4000
$ws.Range("M103").Value2 = -2828
# row 126
$ws.Range("H126").Value2 = 2128.5715
$ws.Range("I126").Value2 = 2040
$ws.Range("J126").Value2 = 2350
$ws.Range("K126").Value2 = 6120
$ws.Range("L126").Value2 = 7050
$ws.Range("M126").Value2 = -3650
$ws.Range("N126").Value2 = -11990

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# row 3
$ws.Range("H3").Value2 = 6187.5
$ws.Range("I3").Value2 = 2333.3333
$ws.Range("J3").Value2 = 8500
$ws.Range("K3").Value2 = 6999.999899999999
$ws.Range("L3").Value2 = 25500
$ws.Range("M3").Value2 = -6887.999899999999
$ws.Range("N3").Value2 = -25724
# row 68
$ws.Range("H68").Value2 = 1500.7142
$ws.Range("J68").Value2 = 1921
$ws.Range("L68").Value2 = 5763
$ws.Range("N68").Value2 = -7385
# row 71
$ws.Range("H71").Value2 = 1500.7142
$ws.Range("J71").Value2 = 1921
$ws.Range("L71").Value2 = 17289
$ws.Range("N71").Value2 = -25401
# row 107
$ws.Range("H107").Value2 = 599509.3
$ws.Range("I107").Value2 = 1693.3334
$ws.Range("J107").Value2 = 778854.1
$ws.Range("K107").Value2 = 5080.0002
$ws.Range("L107").Value2 = 2336562.3
$ws.Range("M107").Value2 = -3160.0002
$ws.Range("N107").Value2 = -2340402.3
# row 131
$ws.Range("H131").Value2 = 25050884
$ws.Range("J131").Value2 = 1147.1428
$ws.Range("L131").Value2 = 3441.4284
$ws.Range("N131").Value2 = -13521.4284
# row 132
$ws.Range("H132").Value2 = 1092.05
$ws.Range("I132").Value2 = 870.5
$ws.Range("J132").Value2 = 1424.375
$ws.Range("K132").Value2 = 7834.5
$ws.Range("L132").Value2 = 12819.375
$ws.Range("M132").Value2 = -5304.5
$ws.Range("N132").Value2 = -17879.375
# row 133
$ws.Range("H133").Value2 = 990
$ws.Range("I133").Value2 = 990
$ws.Range("J133").Value2 = 0
$ws.Range("K133").Value2 = 2970
$ws.Range("L133").Value2 = 0
$ws.Range("M133").Value2 = 2090
$ws.Range("N133").ClearContents()

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# row 126
$ws.Range("H126").Value2 = 3168.5
$ws.Range("I126").Value2 = 3168.5
$ws.Range("K126").Value2 = 9505.5
$ws.Range("M126").Value2 = -7035.5
# row 132
$ws.Range("H132").Value2 = 2843.1707
$ws.Range("I132").Value2 = 2676.9666
$ws.Range("J132").Value2 = 3296.4546
$ws.Range("K132").Value2 = 8030.899800000001
$ws.Range("L132").Value2 = 9889.363799999999
$ws.Range("M132").Value2 = -5500.899800000001
$ws.Range("N132").Value2 = -14949.3638

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# row 25
$ws.Range("H25").Value2 = 9980
$ws.Range("I25").Value2 = 9980
$ws.Range("K25").Value2 = 9980
$ws.Range("M25").Value2 = -9750
# row 40
$ws.Range("H40").Value2 = 1444558.6
$ws.Range("I40").Value2 = 1444558.6
$ws.Range("J40").Value2 = 0
$ws.Range("K40").Value2 = 1444558.6
$ws.Range("L40").Value2 = 0
$ws.Range("M40").Value2 = -1444422.6
$ws.Range("N40").ClearContents()

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# row 2
$ws.Range("H2").Value2 = 5361.25
$ws.Range("I2").Value2 = 4990
$ws.Range("J2").Value2 = 5414.2856
$ws.Range("K2").Value2 = 4990
$ws.Range("L2").Value2 = 5414.2856
$ws.Range("M2").Value2 = -4878
$ws.Range("N2").Value2 = -5638.2856
# row 132
$ws.Range("H132").Value2 = 1003.4583
$ws.Range("I132").Value2 = 1132.2941
$ws.Range("J132").Value2 = 690.5714
$ws.Range("K132").Value2 = 3396.8823
$ws.Range("L132").Value2 = 2071.7142
$ws.Range("M132").Value2 = -866.8823000000002
$ws.Range("N132").Value2 = -7131.7142
